# Add the new "2023" column (S) to the Transport and storage statistics
# table, continuing the existing year series (2006-2022) that already
# occupies columns B-R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing year column (R) into the new
# column (S) so the new cells keep the same borders/number format/fonts
# as the rest of the table (header row 3 plus the 11 data rows 4-14).
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)

# New year header
$ws.Range("S3").Value = 2023

# New year data, one value per indicator row (same order as column R)
$ws.Range("S4").Value = 1604.7
$ws.Range("S5").Value = 1058.2
$ws.Range("S6").Value = 7982
$ws.Range("S7").Value = 7312
$ws.Range("S8").Value = 1638.9
$ws.Range("S9").Value = 432.8
$ws.Range("S10").Value = 154.3
$ws.Range("S11").Value = 625.4
$ws.Range("S12").Value = 153.4
$ws.Range("S13").Value = 975
$ws.Range("S14").Value = 552.4

# Match the selection left behind by the edit (the newly-filled range)
$ws.Range("S4:S14").Select()
